$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at J and K (shifting old J:Q to L:S)
$ws.Range("J1:K1").EntireColumn.Insert()

# Set new headers
$ws.Range("J1").Value = "SOC_start_variance"
$ws.Range("K1").Value = "SOC_end"

# Update changed data values
$ws.Range("E2").Value = 0.5
$ws.Range("E3").Value = 0.5

$ws.Range("H2").Value = 5
$ws.Range("H3").Value = 5

$ws.Range("I2").Value = 0.6
$ws.Range("I3").Value = 0.6

$ws.Range("J2").Value = 0.05
$ws.Range("J3").Value = 0.05

$ws.Range("K2").Value = "0.8,0.9"
$ws.Range("K3").Value = "0.8,0.9"

$ws.Range("M2").Value = 10
$ws.Range("M3").Value = 10

$ws.Range("O2").Value = 0.2
$ws.Range("O3").Value = 0.2

$ws.Range("Q2").Value = 3
$ws.Range("Q3").Value = 3

$ws.Range("R2").Value = 0.3
$ws.Range("R3").Value = 0.3

$ws.Range("J7").Select() | Out-Null
